# Sprint 2 burnup report update
# Updates the "Actual" hours-burned series (column B, rows 9-16) on Sheet1
# with the latest reported progress for the sprint.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B9").Value  = 40.0
$ws.Range("B10").Value = 42.0
$ws.Range("B11").Value = 48.0
$ws.Range("B12").Value = 48.5
$ws.Range("B13").Value = 50.0
$ws.Range("B14").Value = 50.0
$ws.Range("B15").Value = 51.0
$ws.Range("B16").Value = 59.0

# Keep the burnup chart's "plot visible cells only" flag explicit (matches
# the refreshed chart settings saved with this sprint report).
$chartObjects = $ws.ChartObjects()
if ($chartObjects.Count -gt 0) {
    $chart = $chartObjects.Item(1).Chart
    $chart.PlotVisibleOnly = $true
}
